# Live_Config.xlsx update:
#  - Insert a new "CRMInfo" section (header + one Name/Value row) right above
#    the existing "Regex" section, pushing the Regex rows (and everything
#    below) down by two rows.
#  - Resize Table1 / its AutoFilter to cover the two new rows.
#  - Refresh the sheet selection to where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows above the "Regex" section header (row 47) -------
$ws.Rows("47:48").Insert()

# Row-insert clones formatting from the row above (row 46, a Name/Value data
# row), so the new A47 picks up that style and a spurious formatted B47 cell
# appears. Row 48 already ends up with the correct Name/Value formatting
# (style matches row 46), so only A47 / B47 need fixing.

# A47 must look like the bold/italic/underlined section-header cells (e.g.
# the "Regex" header that just got pushed down to A49) instead of a plain
# data cell.
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A47").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# B47 shouldn't exist at all (the header row only has a populated A cell).
$ws.Range("B47").Clear() | Out-Null

# --- Populate the new "CRMInfo" section -------------------------------------
$ws.Range("A47").Value = "CRMInfo"
$ws.Range("A48").Value = "DocumentType"
$ws.Range("B48").Value = "LTRCUST"

# --- Grow the table / autofilter to include the two new rows ---------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:C57")) | Out-Null

# --- Restore a sensible selection near the edited rows ----------------------
$ws.Range("A22").Select() | Out-Null
